$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: 5th fiscal year (31/05/2015) of IPI data for Nike ---
# Copy the visual formatting from column C (dates / numbers, both fully
# numeric columns) so the new column picks up the existing date & number
# styles used throughout the sheet, then overwrite with the new values.

$ws.Range("C1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = 42155

$ws.Range("C2:C16").Copy()
$ws.Range("F2:F16").PasteSpecial(-4122)

$ws.Range("F2").Value = 30601000
$ws.Range("F3").Value = 15885000
$ws.Range("F4").Value = 14716000
$ws.Range("F5").Value = 10541000
$ws.Range("F6").Value = 4175000
$ws.Range("F7").Value = 6000
$ws.Range("F8").Value = -24000
$ws.Range("F9").Value = 4205000
$ws.Range("F10").Value = 932000
$ws.Range("F11").Value = 3273000
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 3273000
$ws.Range("F14").Value = 3922000
$ws.Range("F15").Value = 5079000
$ws.Range("F16").Value = 4824000

# Column F width, matching the author's new column
$ws.Columns.Item(6).ColumnWidth = 10.3

# Selection moved by the author while working further down the sheet
[void]$ws.Range("H22").Select()
